# Add the new participant row (Aldida Chandra Sukma) to the Participants sheet,
# including a mailto: hyperlink on her e-mail address.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data right after the last existing entry (row 14 -> row 15)
$ws.Range("A15").Value = "Aldida"
$ws.Range("B15").Value = "Chandra Sukma"
$ws.Range("C15").Value = "2581065@student.vu.nl"

# Turn the e-mail address into a live hyperlink (adds the built-in
# "Hyperlink" cell style / underlined themed font automatically).
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:2581065@student.vu.nl")

# Match the author's last selection in the sheet.
$ws.Range("C15").Select()
